# Target: the slide whose persisted SlideID is 264 (pc:sldMk cId="3485151880"
# sldId="264" in the change log), shape id=4 ("テキスト ボックス 3", creationId
# {32E015C4-C8C6-3449-B9F4-8ADD6202FBD8}).
#
# Edit: reposition/resize that textbox and rewrite its body copy from
# "Profile:" + <tab> + one long sentence into three lines - "Profile:", then
# the sentence split in two, each continuation line led by an ideographic
# space - while dropping the forced 40pt run size back to the deck's normal
# 18pt.

$p = $ppt.ActivePresentation

$sl = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $cand = $p.Slides.Item($i)
    if ($cand.SlideID -eq 264) { $sl = $cand }
}

$sh = $null
for ($i = 1; $i -le $sl.Shapes.Count; $i++) {
    $cand = $sl.Shapes.Item($i)
    if ($cand.Id -eq 4) { $sh = $cand }
}

$tr = $sh.TextFrame.TextRange

$ideoSpace = [char]0x3000
$line1 = "Profile:"
$line2 = $ideoSpace + "社内システム・社内インフラ全般を"
$line3 = $ideoSpace + "担当しています。"

$tr.Text = $line1 + "`r" + $line2 + "`r" + $line3

# The explicit 40pt run size goes away (falls back to the deck's normal 18pt).
$tr.Font.Size = 18

# Move / resize the shape (EMU -> points; 12700 EMU per point).
$sh.Left   = 156233 / 12700
$sh.Top    = 2436566 / 12700
$sh.Width  = 4108817 / 12700
$sh.Height = 923330 / 12700
